$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename column header labels to use the respective input-file suffixes:
#    "<Name>_old" -> "<Name>_FV2404"   (columns A..J)
#    "<Name>_new" -> "<Name>_FV2410"   (columns L..U, K stays "diff")
$fv2404Names = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)
$fv2410Names = @(
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $fv2404Names.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Names[$i]
}

for ($i = 0; $i -lt $fv2410Names.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $fv2410Names[$i]
}

# 2) Turn the data range into a native Excel Table named "Table1" (adds
#    autoFilter + tableColumns + tableParts wiring).
$range = $ws.Range("A1:U77")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# 3) Freeze the header row (split below row 1, pane anchored at A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
